# Adds a "Discount" (total expense) column to the Inventory sheet and
# starts wiring up the per-item discount/expense totals, while also
# reducing the T.V stock count to reflect Bob's purchase.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inventory")

# Rename the first header from "Name" to "Item" (reuses the existing
# "Item" shared string already used on the Shopping List sheet).
$ws.Range("A1").Value = "Item"

# New header for the discount/expense total column.
$ws.Range("E1").Value = "Discount"

# Per-item total expense (quantity sold * price) aggregated from the
# Shopping List / Expenses sheets.
$ws.Range("E2").Value = 0      # Socks
$ws.Range("E3").Value = 0      # Bananas
$ws.Range("E4").Value = 5      # Ice Cream
$ws.Range("E5").Value = 7.5    # Oranges
$ws.Range("E6").Value = 20     # Nesquik
$ws.Range("E7").Value = 35     # Candy
$ws.Range("E8").Value = 0      # Guitar
$ws.Range("E9").Value = 2      # Paper Plates
$ws.Range("E10").Value = 0     # T.V
$ws.Range("E11").Value = 15    # Laptop

# T.V stock decreases by the one purchased by Bob.
$ws.Range("D10").Value = 1

# Update the sheet's tracked selection to match the cell last worked on.
$ws.Range("F10").Select()
